$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 22:46"

# Rows 47 and 48 had their city name / "Casos activos" values swapped
# (Almeria/72 <-> Lugo/5). Columns B, D, E stay as-is since they are
# identical between the two rows.
$ws.Range("A47").Value = "Lugo"
$ws.Range("C47").Value = 5

$ws.Range("A48").Value = "Almeria"
$ws.Range("C48").Value = 72
